$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Duplicate the last slide ("'Rules' for course") to use as the
#    starting point for the new slide 8 - it already has the right
#    layout (title + single "half" content placeholder).
# ------------------------------------------------------------------
$lastIndex = $p.Slides.Count
$original = $p.Slides.Item($lastIndex)
$original.Duplicate() | Out-Null
$newSlide = $p.Slides.Item($p.Slides.Count)

# ------------------------------------------------------------------
# 2) Fill in the new slide 8 ("Before we begin...")
# ------------------------------------------------------------------
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Before we begin" + [char]0x2026 + " "

$body = $newSlide.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Why are you learning Python?`rUse for research`rImprove chances of a job`rFor fun`rOther reason`rCan you code in another language?`rYes, MATLAB`rYes, R`rYes, other" + [char]0x2026 + "`rNo"

$tr = $body.TextFrame.TextRange
$tr.Paragraphs(2,1).IndentLevel = 2
$tr.Paragraphs(3,1).IndentLevel = 2
$tr.Paragraphs(4,1).IndentLevel = 2
$tr.Paragraphs(5,1).IndentLevel = 2
$tr.Paragraphs(7,1).IndentLevel = 2
$tr.Paragraphs(8,1).IndentLevel = 2
$tr.Paragraphs(9,1).IndentLevel = 2
$tr.Paragraphs(10,1).IndentLevel = 2

$body.Width = 10373140 / 12700

# ------------------------------------------------------------------
# 3) Update the original last slide (now "During the course")
# ------------------------------------------------------------------
$original.Shapes.Item(1).TextFrame.TextRange.Text = "During the course"

$origBody = $original.Shapes.Item(2)
$origBody.TextFrame.TextRange.Text = "If you need to step out at any time, please do so `rAsk questions any time`rBe constructive when helping each other`r`rCoffee and tea in the back`rWater fountain just outside`rToilets outside the door to the right, follow signs`r`r"
